# Regenerate the s_vals sheet data (filtered save games -> new computed stats)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2023-08-18)
$ws.Range("B2").Value = 0.001754667048134761
$ws.Range("C2").Value = 0.000002220651329265522
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 247.6997404328161

# Row 3 (2023-08-04)
$ws.Range("B3").Value = 0.3464964993005633
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 10.24810883928323

# Row 4 (2023-08-03)
$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 5.488907176552729

# Row 5 (2023-05-15) - note C5/G5 are huge values, written without
# scientific-notation syntax (parser only accepts plain decimals);
# trailing ".0" keeps them parsed as floating point instead of
# overflowing a 64-bit integer literal.
$ws.Range("B5").Value = 0.3464964993005633
$ws.Range("C5").Value = 31858158836671370000.0
$ws.Range("D5").Value = 0.7127328510149897
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("G5").Value = 31858158836671370000.0

# Row 6 (2023-03-11)
$ws.Range("B6").Value = 1.505614041169197
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.1529057820181812
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("G6").Value = 3.811642989160245

# Row 7 (2023-03-03)
$ws.Range("B7").Value = 3.182878228561681
$ws.Range("C7").Value = 0.3375848360084654
$ws.Range("D7").Value = 0.1529057820181812
$ws.Range("E7").Value = 0.4998867070740569
$ws.Range("G7").Value = 4.173255553662385
